# chore: update Sheets via scheduled runner
# Refresh leve-crafting price/profit figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the ALC/ARM/BSM/CRP/GSM/LTW/WVR sheets with the latest
# market-board snapshot values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 98 (hunk 1, G=36237)
$ws.Range("H98").Value = 1485
$ws.Range("I98").Value = 1485
$ws.Range("K98").Value = 1485
$ws.Range("M98").Value = 13

# Row 107 (hunk 2, G=27766)
$ws.Range("H107").Value = 338.2381
$ws.Range("I107").Value = 305.15
$ws.Range("K107").Value = 305.15
$ws.Range("M107").Value = 1614.85

# Row 122 (hunk 3, G=36237)
$ws.Range("H122").Value = 1485
$ws.Range("I122").Value = 1485
$ws.Range("K122").Value = 4455
$ws.Range("M122").Value = -2005

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (hunk 4, G=27714)
$ws.Range("H45").Value = 1160.8889
$ws.Range("I45").Value = 1078.2858
$ws.Range("J45").Value = 1450
$ws.Range("K45").Value = 1078.2858
$ws.Range("L45").Value = 1450
$ws.Range("M45").Value = -701.2858000000001
$ws.Range("N45").Value = -2204

# Row 94 (hunk 5, G=18055)
$ws.Range("H94").Value = 35450
$ws.Range("J94").Value = 35450
$ws.Range("L94").Value = 35450
$ws.Range("N94").Value = -37252

# Row 97 (hunk 6, G=19941)
$ws.Range("H97").Value = 1829.4546
$ws.Range("I97").Value = 664.5833
$ws.Range("J97").Value = 4935.778
$ws.Range("K97").Value = 664.5833
$ws.Range("L97").Value = 4935.778
$ws.Range("M97").Value = -168.5833
$ws.Range("N97").Value = -5927.778

# Row 102 (hunk 7, G=19945)
$ws.Range("H102").Value = 1505.0834
$ws.Range("I102").Value = 1218.5714
$ws.Range("J102").Value = 1906.2
$ws.Range("K102").Value = 1218.5714
$ws.Range("L102").Value = 1906.2
$ws.Range("M102").Value = 403.4286
$ws.Range("N102").Value = -5150.2

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (hunk 8, G=19943)
$ws.Range("H99").Value = 896.9524
$ws.Range("I99").Value = 777.94116
$ws.Range("K99").Value = 777.94116
$ws.Range("M99").Value = 720.05884

# Row 105 (hunk 9, G=19947)
$ws.Range("H105").Value = 2975.5518
$ws.Range("I105").Value = 2978.4614
$ws.Range("J105").Value = 2973.1875
$ws.Range("K105").Value = 2978.4614
$ws.Range("L105").Value = 2973.1875
$ws.Range("M105").Value = -1231.4614
$ws.Range("N105").Value = -6467.1875

# Row 122 (hunk 10, G=34096)
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").Value = ""

# Row 125 (hunk 11, G=34235)
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = ""
$ws.Range("N125").Value = ""

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (hunk 12, G=27691)
$ws.Range("H16").Value = 2659.1765
$ws.Range("I16").Value = 1305.7273
$ws.Range("J16").Value = 5140.5
$ws.Range("K16").Value = 1305.7273
$ws.Range("L16").Value = 5140.5
$ws.Range("M16").Value = -1018.7273
$ws.Range("N16").Value = -5714.5

# Row 28 (hunk 13, G=18348)
$ws.Range("H28").Value = 47821.5
$ws.Range("J28").Value = 47821.5
$ws.Range("L28").Value = 47821.5
$ws.Range("N28").Value = -48311.5

# Row 107 (hunk 14, G=27689)
$ws.Range("H107").Value = 811.64703
$ws.Range("I107").Value = 860.5714
$ws.Range("J107").Value = 583.3333
$ws.Range("K107").Value = 860.5714
$ws.Range("L107").Value = 583.3333
$ws.Range("M107").Value = 1059.4286
$ws.Range("N107").Value = -4423.3333

# Row 113 (hunk 15, G=27691)
$ws.Range("H113").Value = 2659.1765
$ws.Range("I113").Value = 1305.7273
$ws.Range("J113").Value = 5140.5
$ws.Range("K113").Value = 1305.7273
$ws.Range("L113").Value = 5140.5
$ws.Range("M113").Value = 864.2727
$ws.Range("N113").Value = -9480.5

# Row 132 (hunk 16, G=44019)
$ws.Range("H132").Value = 71439300
$ws.Range("I132").Value = 250029900
$ws.Range("J132").Value = 3053
$ws.Range("K132").Value = 750089700
$ws.Range("L132").Value = 9159
$ws.Range("M132").Value = -750087170
$ws.Range("N132").Value = -14219

# Row 134 (hunk 17, G=44020)
$ws.Range("H134").Value = 2055.4167
$ws.Range("I134").Value = 2003.5294
$ws.Range("J134").Value = 2181.4285
$ws.Range("K134").Value = 6010.5882
$ws.Range("L134").Value = 6544.2855
$ws.Range("M134").Value = -3475.5882
$ws.Range("N134").Value = -11614.2855

$ws = $wb.Worksheets.Item("GSM")
# Row 34 (hunk 18, G=10924)
$ws.Range("H34").Value = 17036.5
$ws.Range("J34").Value = 17036.5
$ws.Range("L34").Value = 17036.5
$ws.Range("N34").Value = -17572.5

# Row 76 (hunk 19, G=10924)
$ws.Range("H76").Value = 17036.5
$ws.Range("J76").Value = 17036.5
$ws.Range("L76").Value = 17036.5
$ws.Range("N76").Value = -17666.5

# Row 79 (hunk 20, G=10924)
$ws.Range("H79").Value = 17036.5
$ws.Range("J79").Value = 17036.5
$ws.Range("L79").Value = 17036.5
$ws.Range("N79").Value = -19220.5

# Row 126 (hunk 21, G=36184)
$ws.Range("H126").Value = 1700
$ws.Range("I126").Value = 1600
$ws.Range("J126").Value = 1800
$ws.Range("K126").Value = 4800
$ws.Range("L126").Value = 5400
$ws.Range("M126").Value = -2330
$ws.Range("N126").Value = -10340

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (hunk 22, G=36249)
$ws.Range("H7").Value = 1978.1538
$ws.Range("I7").Value = 1346.2222
$ws.Range("J7").Value = 3400
$ws.Range("K7").Value = 1346.2222
$ws.Range("L7").Value = 3400
$ws.Range("M7").Value = -1234.2222
$ws.Range("N7").Value = -3624

# Row 16 (hunk 23, G=5289)
$ws.Range("H16").Value = 457.72726
$ws.Range("I16").Value = 465.8889
$ws.Range("J16").Value = 421
$ws.Range("K16").Value = 465.8889
$ws.Range("L16").Value = 421
$ws.Range("M16").Value = -295.8889
$ws.Range("N16").Value = -761

# Row 126 (hunk 24, G=36249)
$ws.Range("H126").Value = 1978.1538
$ws.Range("I126").Value = 1346.2222
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 4038.6666
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -1568.6666
$ws.Range("N126").Value = -15140

# Row 132 (hunk 25, G=44058)
$ws.Range("H132").Value = 8070.3335
$ws.Range("I132").Value = 9860.105
$ws.Range("J132").Value = 3819.625
$ws.Range("K132").Value = 29580.315
$ws.Range("L132").Value = 11458.875
$ws.Range("M132").Value = -27050.315
$ws.Range("N132").Value = -16518.875

$ws = $wb.Worksheets.Item("WVR")
# Row 22 (hunk 26, G=3041)
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = ""

# Row 96 (hunk 27, G=19977)
$ws.Range("H96").Value = 872.5833
$ws.Range("I96").Value = 720.875
$ws.Range("J96").Value = 1176
$ws.Range("K96").Value = 720.875
$ws.Range("L96").Value = 1176
$ws.Range("M96").Value = 652.125
$ws.Range("N96").Value = -3922

# Row 103 (hunk 28, G=18548)
$ws.Range("H103").Value = 15500
$ws.Range("J103").Value = 23000
$ws.Range("L103").Value = 23000
$ws.Range("N103").Value = -25344

# Row 132 (hunk 29, G=44029)
$ws.Range("H132").Value = 7612.28
$ws.Range("I132").Value = 9289.117
$ws.Range("J132").Value = 4049
$ws.Range("K132").Value = 27867.351
$ws.Range("L132").Value = 12147
$ws.Range("M132").Value = -25337.351
$ws.Range("N132").Value = -17207
